$d = $word.ActiveDocument

# Insert a new run with text "me" at the very start of the first paragraph
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertBefore("me")

# Move the "_GoBack" bookmark from its old location (an otherwise-empty
# paragraph between the two tables) to right after the newly inserted
# "me" run, at the start of the first paragraph.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$newSpot = $d.Range(2, 2)
$d.Bookmarks.Add("_GoBack", $newSpot)
